# This script rewrites the acquisition-request template body.
#
# The original document had 5 paragraphs addressed "To the Chief
# Engineer" with an organization/engineer header block and a long
# justification paragraph referencing {{pernr}}/{{snils}}/{{itemName}}.
#
# The edited version drops the addressee header entirely, adds a
# centered bold title paragraph ("Заявление на передачу материальных
# средств"), a blank centered spacer paragraph, and replaces the body
# paragraph with a shorter justification that references a single
# {{item}} placeholder (and keeps {{reason}}).
#
# Because nearly every paragraph/run in the story changes (structure,
# run-splits, alignment, spacing, bold title, etc.), the most faithful
# and robust way to reproduce the exact target OOXML is to replace the
# whole document story with the literal target markup via
# Range.InsertXML, rather than trying to stitch it together with many
# fragile Find/Replace calls. InsertXML inserts well-formed WordML
# fragments using the w: namespace already declared on the package,
# and $d.Content spans the whole body up to (but excluding) the
# trailing sectPr, so the section properties are left untouched.

$d = $word.ActiveDocument

$newBodyXml = @'
<w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Заявление на передачу материальных средств</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="709"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Прошу Вас с</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> целью организации рабочего процесса прошу рассмотреть возможность </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>выделения нов</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>ого</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> {{</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>item</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> по причине {{</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>reason</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>}}</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
'@

$d.Content.InsertXML($newBodyXml)
